# Update the "Metadata" worksheet of the StructureDefinition workbook:
#   - Version bumped from 5.0.0 to 6.0.0
#   - Date bumped to the new publication timestamp
#   - Publisher value filled in ("Alvearie Team")
#   - The duplicated "Contact" / "No display for ContactDetail" row is
#     replaced by a single "Jurisdiction" / "United States of America" row
#     (net effect: one fewer row, so everything below shifts up by one)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "6.0.0"
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$ws.Range("B9").Value = "Alvearie Team"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# The old row 11 ("Contact" | "No display for ContactDetail") is removed
# entirely; Excel shifts the remaining rows (Description, Purpose, ...) up.
$ws.Range("A11").EntireRow.Delete()
